# =========================================================================
# Edit: add "2022-Q1" fund-holdings sheet (inserted before the existing
# "总计" summary sheet) and prepend the corresponding summary row to
# "总计".
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right after "2021-Q4" (i.e. right
#    before "总计", which is currently the last sheet).
# -------------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $src)
$q1.Name = "2022-Q1"

# Clone number/border/alignment formatting from the template sheet so the
# new sheet's header row + index column look like every other quarter
# sheet (bold, centered, thin border - style index 2 in this workbook).
$src.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$q1.Range("A2:A33").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund holdings, one "|"-delimited record per row:
#   index|基金代码|基金名称|基金规模|股票总仓位|仓位占比|持有市值(亿元)|仓位排名
$data = @"
0|570001|诺德价值优势混合|40.24|92.31|8.11|3.2635|4
1|160926|大成创业板两年定期开放混合A|40.93|64.09|4.57|1.8705|3
2|377020|上投摩根内需动力混合|21.62|90.74|6.23|1.3469|3
3|570008|诺德周期策略混合|11.18|91.96|8.48|0.9481|2
4|501079|大成科创主题 3 年封闭运作灵活配置混合|17.69|79.13|4.67|0.8261|3
5|070099|嘉实优质企业混合|22.02|91.84|3.51|0.7729|10
6|012150|诺德价值发现一年持有期混合型证券投资基金|9.42|91.48|8.10|0.7630|2
7|399011|中海医疗保健主题股票|10.58|87.56|3.64|0.3851|10
8|009774|财通资管优选回报一年持有期混合|8.48|93.17|3.86|0.3273|8
9|001766|上投摩根医疗健康股票|10.35|80.54|2.79|0.2888|7
10|009798|大成创业板两年定期开放混合C|5.67|64.09|4.57|0.2591|3
11|159883|永赢中证全指医疗器械交易型开放式指数证券投资基金|9.17|99.23|2.76|0.2531|8
12|010054|万家健康产业混合A|8.13|86.63|2.94|0.2390|9
13|570005|诺德成长优势混合|4.30|62.89|4.89|0.2103|3
14|010371|大成成长进取混合A|5.55|80.17|3.66|0.2031|9
15|000073|上投摩根成长动力混合|3.03|90.81|6.35|0.1924|3
16|005682|财通资管消费精选灵活配置混合A|3.64|94.77|4.58|0.1667|7
17|000870|嘉实新收益灵活配置混合|4.17|91.30|3.53|0.1472|10
18|012036|诺德兴远优选一年持有期混合型证券投资基金|2.75|52.19|4.92|0.1353|2
19|001192|上投摩根整合驱动灵活配置混合|3.96|76.32|3.17|0.1255|10
20|010055|万家健康产业混合C|3.36|86.63|2.94|0.0988|9
21|010372|大成成长进取混合C|1.71|80.17|3.66|0.0626|9
22|000326|南方中小盘成长股票|2.96|90.75|1.92|0.0568|5
23|159898|招商中证全指医疗器械交易型开放式指数证券投资基金|1.61|99.41|2.74|0.0441|9
24|008277|财通资管行业精选混合|1.15|89.68|3.40|0.0391|7
25|159873|天弘中证全指医疗保健设备与服务ETF|1.43|99.59|2.40|0.0343|10
26|159891|建信中证全指医疗保健设备与服务交易型开放式指数证券投资基金|1.40|95.24|2.32|0.0325|10
27|003561|诺德成长精选灵活配置混合A|0.60|53.97|4.90|0.0294|2
28|006235|东方城镇消费主题混合|0.50|90.32|4.36|0.0218|9
29|516610|大成中证全指医疗保健设备与服务交易型开放式指数证券投资基金|0.69|96.20|2.35|0.0162|10
30|011020|财通资管消费精选灵活配置混合C|0.14|94.77|4.58|0.0064|7
31|003562|诺德成长精选灵活配置混合C|0.00|53.97|4.90|0|2
"@

$lines = $data -split "`n"
$r = 2
foreach ($line in $lines) {
  $line = $line.Trim()
  if ($line.Length -eq 0) { continue }
  $p = $line -split "\|"

  $q1.Cells.Item($r, 1).Value = [int]$p[0]

  $q1.Cells.Item($r, 2).Value = "'" + $p[1]
  $q1.Cells.Item($r, 3).Value = $p[2]
  $q1.Cells.Item($r, 4).Value = "'" + $p[3]
  $q1.Cells.Item($r, 5).Value = "'" + $p[4]
  $q1.Cells.Item($r, 6).Value = "'" + $p[5]

  # Every holding row stores "持有市值(亿元)" as text EXCEPT the very last
  # (smallest) holding, whose value of exactly "0" is stored as a real
  # number in the source workbook.
  if ($p[6] -eq "0") {
    $q1.Cells.Item($r, 7).Value = 0
  } else {
    $q1.Cells.Item($r, 7).Value = "'" + $p[6]
  }

  $q1.Cells.Item($r, 8).Value = [int]$p[7]

  $r++
}

# The fund-code/name/size/position columns (B:G) were entered with a
# leading "'" to force text storage (preserving leading zeros in fund
# codes and exact decimal text like "40.24"); strip the resulting
# quote-prefix formatting so the cells end up with no explicit style,
# matching the rest of the workbook.
$q1.Range("B2:G33").ClearFormats()

# -------------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" (grand total) summary sheet.
# -------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")
$tot.Rows.Item(2).Insert()

$tot.Range("B2:D2").ClearFormats()
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 32
$tot.Range("D2").Value = 13.17

$tot.Range("A3").Copy()
$tot.Range("A2").PasteSpecial(-4122)
$tot.Range("A2").Value = 0

# Renumber the index column (A) for the rows that shifted down one slot.
$tot.Range("A3").Value = 1
$tot.Range("A4").Value = 2
$tot.Range("A5").Value = 3
$tot.Range("A6").Value = 4
$tot.Range("A7").Value = 5
